$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Cadastro")
$ws2 = $wb.Worksheets.Item("Produtos")

# --- Username value update (shared by B2 and B14 on "Cadastro") ---
$ws1.Range("B2").Value = "john1475"
$ws1.Range("B14").Value = "john1475"

# --- New style A: numFmtId=1 (format "0"), no border - applied to blank cell C2 ---
$ws1.Range("C2").NumberFormat = "0"

# --- New style B: numFmtId=1, bordered, quote-prefixed (text "999"), right aligned - D5 on "Produtos" ---
$ws2.Range("D5").Value = "'999"
$ws2.Range("D5").NumberFormat = "0"
$ws2.Range("D5").HorizontalAlignment = -4152

# --- Selection / active sheet changes ---
$ws1.Range("C2").Select()
$ws2.Activate()
$ws2.Range("C7").Select()
